# Apply updated Betfair Back/Lay odds values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = 2.32
$ws.Range("T2").Value = 1.57
$ws.Range("U2").Value = 2.38
$ws.Range("AJ2").Value = 21
$ws.Range("AN2").Value = 8.199999999999999

# Row 3
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 5
$ws.Range("N3").Value = 5.4
$ws.Range("R3").Value = 1.59
$ws.Range("T3").Value = 1.88
$ws.Range("U3").Value = 2.06
$ws.Range("W3").Value = 3.1
$ws.Range("AN3").Value = 5.6

# Row 4
$ws.Range("G4").Value = 2.08
$ws.Range("H4").Value = 3.65
$ws.Range("L4").Value = 1.31
$ws.Range("P4").Value = 2.42
$ws.Range("Q4").Value = 1.66
$ws.Range("R4").Value = 1.57
$ws.Range("AJ4").Value = 25
$ws.Range("AO4").Value = 27

# Row 5
$ws.Range("F5").Value = 1.75
$ws.Range("G5").Value = 1.97
$ws.Range("H5").Value = 3.85
$ws.Range("I5").Value = 5
$ws.Range("K5").Value = 5.2
$ws.Range("L5").Value = 1.22
$ws.Range("N5").Value = 3.45
$ws.Range("O5").Value = 1.13
$ws.Range("P5").Value = 2.46
$ws.Range("Q5").Value = 1.13
$ws.Range("R5").Value = 1.7
$ws.Range("S5").Value = 1.89
$ws.Range("T5").Value = 1.33
$ws.Range("U5").Value = 2.28
$ws.Range("V5").Value = 1.25
$ws.Range("W5").Value = 2.02
$ws.Range("X5").Value = 44
$ws.Range("Y5").Value = 34
$ws.Range("AB5").Value = 21
$ws.Range("AC5").Value = 15
$ws.Range("AD5").Value = 22
$ws.Range("AE5").Value = 48
$ws.Range("AF5").Value = 21
$ws.Range("AG5").Value = 14
$ws.Range("AH5").Value = 19
$ws.Range("AI5").Value = 44
$ws.Range("AJ5").Value = 27
$ws.Range("AK5").Value = 20
$ws.Range("AL5").Value = 24
$ws.Range("AM5").Value = 60
$ws.Range("AN5").Value = 7.6
$ws.Range("AO5").Value = 27

# Row 6
$ws.Range("I6").Value = 2.98
$ws.Range("J6").Value = 2.72

# Row 8
$ws.Range("S8").Value = 4
$ws.Range("V8").Value = 1.42

# Row 9
$ws.Range("F9").Value = 2.64
$ws.Range("G9").Value = 2.66
$ws.Range("H9").Value = 2.64
$ws.Range("I9").Value = 2.68
$ws.Range("O9").Value = 1.16
$ws.Range("R9").Value = 1.74
$ws.Range("S9").Value = 2.28
$ws.Range("T9").Value = 1.49
$ws.Range("V9").Value = 1.59
$ws.Range("W9").Value = 1.6
$ws.Range("AI9").Value = 27
$ws.Range("AO9").Value = 13

# Row 10
$ws.Range("S10").Value = 2.6
$ws.Range("T10").Value = 1.98
$ws.Range("U10").Value = 1.96
$ws.Range("AO10").Value = 5

# Row 11
$ws.Range("Q11").Value = 1.48
$ws.Range("T11").Value = 2.34
$ws.Range("U11").Value = 1.7
$ws.Range("AH11").Value = 160
$ws.Range("AN11").Value = 3.3

# Row 12
$ws.Range("K12").Value = 7.2
$ws.Range("N12").Value = 8.6
$ws.Range("P12").Value = 3.5
$ws.Range("S12").Value = 1.92
$ws.Range("X12").Value = 42
$ws.Range("AM12").Value = 95

# Row 13
$ws.Range("G13").Value = 6.6
$ws.Range("Q13").Value = 1.82
$ws.Range("AL13").Value = 80

# Row 14
$ws.Range("N14").Value = 4.9
$ws.Range("P14").Value = 2.3
$ws.Range("R14").Value = 1.52

# Row 15
$ws.Range("I15").Value = 3.15
$ws.Range("V15").Value = 1.46

# Row 16
$ws.Range("F16").Value = 2.48
$ws.Range("G16").Value = 2.82
$ws.Range("I16").Value = 3.25
$ws.Range("J16").Value = 3.15
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 3.3
$ws.Range("P16").Value = 1.8
$ws.Range("Q16").Value = 2.02
$ws.Range("S16").Value = 3.55
$ws.Range("T16").Value = 1.76
$ws.Range("W16").Value = 1.57
$ws.Range("Z16").Value = 25
$ws.Range("AA16").Value = 60
$ws.Range("AB16").Value = 13
$ws.Range("AE16").Value = 44
$ws.Range("AH16").Value = 21
$ws.Range("AI16").Value = 60
$ws.Range("AJ16").Value = 48
$ws.Range("AM16").Value = 120
$ws.Range("AO16").Value = 42
